$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update word-search result strings whose coordinates were recalculated (0-based -> 1-based)
$ws.Range("D2").Value = "['Palabra AIRE encontrada entre (1, 1) y (4, 4).', 'Palabra AGUA encontrada entre (1, 1) y (1, 4).', 'Palabra TIERRA no encontrada.', 'Palabra NIU encontrada entre (3, 1) y (1, 3).']"
$ws.Range("D14").Value = "['Palabra AIRE encontrada entre (1, 1) y (4, 4).', 'Palabra AGUA encontrada entre (1, 1) y (1, 4).', 'Palabra TIERRA no encontrada.', 'Palabra NIU encontrada entre (3, 1) y (1, 3).']"
$ws.Range("D13").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D18").Value = "['Palabra AIRE encontrada entre (1, 1) y (4, 4).', 'Palabra AGUA encontrada entre (1, 1) y (1, 4).', 'Palabra TIERRA no encontrada.', 'Palabra NIU encontrada entre (3, 1) y (1, 3).']"

# 2) Insert three new rows for the new gr_as_013/014/015 tests, copying formatting from row 22
$ws.Rows("23:25").Insert()
$ws.Range("A22:E22").Copy()
$ws.Range("A23:E25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Populate the new rows
$ws.Range("A23").Value = "gr_as_013"
$ws.Range("B23").Value = "YES"
$ws.Range("C23").Value = "word out of bounds horizontally"
$ws.Range("D23").Value = "['Palabra PAJARO no encontrada.', 'Palabra PASA encontrada entre (1, 1) y (1, 4).']"
$ws.Range("E23").Value = "Word PAJARO couldn't be found as it is splitted (see test file)"

$ws.Range("A24").Value = "gr_as_014"
$ws.Range("B24").Value = "YES"
$ws.Range("C24").Value = "word out of bounds vertically"
$ws.Range("D24").Value = "['Palabra PAJARO no encontrada.', 'Palabra PASA encontrada entre (1, 1) y (1, 4).']"
$ws.Range("E24").Value = "Word PAJARO couldn't be found as it is splitted (see test file)"

$ws.Range("A25").Value = "gr_as_015"
$ws.Range("B25").Value = "YES"
$ws.Range("C25").Value = "word out of bounds vertically"
$ws.Range("D25").Value = "['Palabra PAJARO no encontrada.', 'Palabra PASA encontrada entre (1, 1) y (1, 4).']"
$ws.Range("E25").Value = "Word PAJARO couldn't be found as it gets out of the table"

# 4) Update shifted ms_as rows (now 26-29) with their new content
$ws.Range("D26").Value = "['Palabra ABADESA encontrada entre (9, 7) y (3, 1).', 'Palabra BALANO encontrada entre (7, 2) y (7, 7).', 'Palabra BATERIA encontrada entre (8, 10) y (2, 10).', 'Palabra BORDADO encontrada entre (4, 13) y (10, 13).', 'Palabra CIERVA encontrada entre (4, 3) y (4, 8).', 'Palabra INTERNO encontrada entre (10, 4) y (10, 10).', 'Palabra LLUVIA encontrada entre (4, 1) y (9, 1).', 'Palabra MARTY encontrada entre (6, 11) y (2, 11).', 'Palabra MINIMO encontrada entre (1, 4) y (1, 9).', 'Palabra PAJIZO encontrada entre (3, 7) y (3, 2).', 'Palabra PAÑUELO encontrada entre (2, 9) y (8, 9).', 'Palabra SCHOTTE encontrada entre (8, 12) y (2, 12).', 'Palabra SOFIA encontrada entre (5, 4) y (5, 8).', 'Palabra SOLARIS encontrada entre (11, 5) y (11, 11).']"
$ws.Range("D27").Value = "['Palabra ADOPCION encontrada entre (22, 2) y (22, 9).', 'Palabra ADORNO encontrada entre (8, 4) y (13, 9).', 'Palabra ALHAJA encontrada entre (23, 1) y (28, 6).', 'Palabra ALIANZA encontrada entre (6, 9) y (12, 9).', 'Palabra AMOR encontrada entre (18, 6) y (18, 9).', 'Palabra ANILLO encontrada entre (10, 1) y (15, 6).', 'Palabra BONDAD encontrada entre (1, 8) y (6, 3).', 'Palabra CALMA encontrada entre (18, 5) y (14, 1).', 'Palabra CANDOR encontrada entre (20, 13) y (25, 8).', 'Palabra CARIÑOS encontrada entre (13, 3) y (19, 9).', 'Palabra DIADEMA encontrada entre (4, 2) y (10, 2).', 'Palabra ESTILO encontrada entre (10, 3) y (5, 8).', 'Palabra ETICA encontrada entre (2, 9) y (6, 9).', 'Palabra FIESTA encontrada entre (23, 8) y (28, 3).', 'Palabra FORTUNA encontrada entre (22, 13) y (28, 13).', 'Palabra GEMA encontrada entre (20, 7) y (20, 10).', 'Palabra GOZO encontrada entre (20, 6) y (20, 3).', 'Palabra HUMILDAD encontrada entre (11, 3) y (18, 10).', 'Palabra LEALTAD encontrada entre (27, 6) y (27, 12).', 'Palabra MAJESTAD encontrada entre (12, 11) y (19, 11).', 'Palabra MESURA encontrada entre (4, 11) y (9, 11).', 'Palabra MODESTIA encontrada entre (21, 1) y (21, 8).', 'Palabra MUSICA encontrada entre (12, 11) y (17, 6).', 'Palabra NACAR encontrada entre (28, 5) y (28, 9).', 'Palabra PERDON encontrada entre (3, 10) y (8, 10).', 'Palabra PLATINO encontrada entre (2, 12) y (8, 12).', 'Palabra PREMIO encontrada entre (15, 10) y (20, 5).', 'Palabra PUREZA encontrada entre (9, 3) y (14, 8).', 'Palabra RECREO encontrada entre (3, 11) y (8, 6).', 'Palabra RISA encontrada entre (15, 1) y (18, 4).', 'Palabra SANTIDAD encontrada entre (19, 13) y (26, 6).', 'Palabra SORTIJA encontrada entre (9, 13) y (3, 13).', 'Palabra TACTO encontrada entre (2, 1) y (6, 1).', 'Palabra TEATRO encontrada entre (26, 7) y (26, 12).', 'Palabra VIRTUD encontrada entre (1, 9) y (6, 4).']"
$ws.Range("B28").Value = "NO"
$ws.Range("C28").Value = "Bad test 1"

# 5) Refresh the view state (selection / scroll position) to match the edited document
$ws.Range("E14").Select()
